# Splitting the Weather data to each location for both models - base and intraday
# Shifts the data date from 23.07.2024 to 30.07.2024 (+7 days) and updates
# the Prediction values for rows 29-84 to the new model output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift column A (Data/date-time) forward by 7 days for rows 2..96 ---
for ($r = 2; $r -le 96; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value2
    $cell.Value = $oldVal + 7
}

# --- 2. Update column D (Lookup) text: "23.07.2024" -> "30.07.2024" ---
for ($r = 2; $r -le 96; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $text = [string]$cell.Value2
    $newText = $text.Replace("23.07.2024", "30.07.2024")
    $cell.Value = $newText
}

# --- 3. Update column C (Prediction) values for rows 29-84 ---
$predictionUpdates = @{
    29 = 0.016
    30 = 0.031
    31 = 0.053
    32 = 0.083
    33 = 0.121
    34 = 0.165
    35 = 0.211
    36 = 0.258
    37 = 0.307
    38 = 0.354
    39 = 0.4
    40 = 0.443
    41 = 0.483
    42 = 0.522
    43 = 0.558
    44 = 0.585
    45 = 0.629
    46 = 0.662
    47 = 0.694
    48 = 0.694
    49 = 0.712
    50 = 0.725
    51 = 0.725
    52 = 0.725
    53 = 0.725
    54 = 0.725
    55 = 0.723
    56 = 0.723
    57 = 0.72
    58 = 0.713
    59 = 0.713
    60 = 0.7
    61 = 0.675
    62 = 0.667
    63 = 0.664
    64 = 0.655
    65 = 0.633
    66 = 0.604
    67 = 0.572
    68 = 0.545
    69 = 0.506
    70 = 0.464
    71 = 0.403
    72 = 0.344
    73 = 0.312
    74 = 0.275
    75 = 0.23
    76 = 0.181
    77 = 0.148
    78 = 0.119
    79 = 0.091
    80 = 0.073
    81 = 0.061
    82 = 0.048
    83 = 0.04
    84 = 0.034
}

foreach ($row in $predictionUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $predictionUpdates[$row]
}
